$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price/volume snapshot. Column D prices are
# stored as text (values such as "26.548.38" or "1.011" are not valid
# numbers), so a leading apostrophe forces Excel to keep them as text
# instead of auto-converting them to numbers.
# Rows 12/13, 36/37, 39/40 and 47/48 also swap their Coin name/Link
# (columns B/C) as the ranking order changed.

# Row 2
$ws.Range("D2").Value = "'26.548.38"
$ws.Range("E2").Value = "  -2.94%  "

# Row 3
$ws.Range("D3").Value = "'1.808.19"
$ws.Range("E3").Value = "  -2.54%  "

# Row 4
$ws.Range("E4").Value = "  +0.85%  "

# Row 5
$ws.Range("D5").Value = "'1.011"
$ws.Range("E5").Value = "  +0.82%  "

# Row 6
$ws.Range("D6").Value = "'308.97"
$ws.Range("E6").Value = "  -1.70%  "

# Row 7
$ws.Range("D7").Value = "'0.4534"
$ws.Range("E7").Value = "  -1.77%  "

# Row 8
$ws.Range("D8").Value = "'0.3667"
$ws.Range("E8").Value = "  -1.32%  "

# Row 9
$ws.Range("D9").Value = "'0.07101"
$ws.Range("E9").Value = "  -2.87%  "

# Row 10
$ws.Range("D10").Value = "'0.8685"
$ws.Range("E10").Value = "  -1.94%  "

# Row 11
$ws.Range("D11").Value = "'0.07822"
$ws.Range("E11").Value = "  +0.13%  "

# Row 12
$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "'19.19"
$ws.Range("E12").Value = "  -3.70%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "'1.855.29"
$ws.Range("E13").Value = "  +0.25%  "

# Row 14
$ws.Range("E14").Value = "  -1.96%  "

# Row 15
$ws.Range("D15").Value = "'6.308"
$ws.Range("E15").Value = "  -3.79%  "

# Row 16
$ws.Range("D16").Value = "'86.34"
$ws.Range("E16").Value = "  -5.93%  "

# Row 17
$ws.Range("D17").Value = "'1.013"
$ws.Range("E17").Value = "  +0.97%  "

# Row 18
$ws.Range("D18").Value = "'0.000008580"
$ws.Range("E18").Value = "  -4.22%  "

# Row 19
$ws.Range("D19").Value = "'1.011"
$ws.Range("E19").Value = "  +0.74%  "

# Row 20
$ws.Range("D20").Value = "'26.597.23"
$ws.Range("E20").Value = "  -2.83%  "

# Row 21
$ws.Range("D21").Value = "'14.29"
$ws.Range("E21").Value = "  -3.21%  "

# Row 22
$ws.Range("D22").Value = "'4.951"
$ws.Range("E22").Value = "  -3.23%  "

# Row 23
$ws.Range("D23").Value = "'2.081.05"
$ws.Range("E23").Value = "  +0.20%  "

# Row 24
$ws.Range("E24").Value = "  -1.73%  "

# Row 25
$ws.Range("E25").Value = "  +3.01%  "

# Row 26
$ws.Range("D26").Value = "'151.12"
$ws.Range("E26").Value = "  -0.37%  "

# Row 27
$ws.Range("D27").Value = "'17.90"
$ws.Range("E27").Value = "  -2.93%  "

# Row 28
$ws.Range("D28").Value = "'1.998"
$ws.Range("E28").Value = "  -2.49%  "

# Row 29
$ws.Range("D29").Value = "'112.90"
$ws.Range("E29").Value = "  -2.72%  "

# Row 30
$ws.Range("D30").Value = "'4.863"
$ws.Range("E30").Value = "  -4.48%  "

# Row 31
$ws.Range("D31").Value = "'0.08684"
$ws.Range("E31").Value = "  -1.83%  "

# Row 32
$ws.Range("D32").Value = "'3.040"
$ws.Range("E32").Value = "  -1.68%  "

# Row 33
$ws.Range("D33").Value = "'0.7315"
$ws.Range("E33").Value = "  -5.33%  "

# Row 34
$ws.Range("D34").Value = "'4.432"
$ws.Range("E34").Value = "  -1.64%  "

# Row 35
$ws.Range("D35").Value = "'1.107"
$ws.Range("E35").Value = "  -5.85%  "

# Row 36
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D36").Value = "'2.470"
$ws.Range("E36").Value = "  -7.96%  "

# Row 37
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").Value = "'1.079"
$ws.Range("E37").Value = "  -0.05%  "

# Row 38
$ws.Range("D38").Value = "'0.01913"
$ws.Range("E38").Value = "  -2.27%  "

# Row 39
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.05093"
$ws.Range("E39").Value = "  -2.71%  "

# Row 40
$ws.Range("B40").Value = "MXToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D40").Value = "'2.872"
$ws.Range("E40").Value = "  -2.61%  "

# Row 41
$ws.Range("D41").Value = "'6.893"
$ws.Range("E41").Value = "  -2.29%  "

# Row 42
$ws.Range("D42").Value = "'0.4906"
$ws.Range("E42").Value = "  -4.55%  "

# Row 43
$ws.Range("D43").Value = "'0.1570"
$ws.Range("E43").Value = "  -4.05%  "

# Row 44
$ws.Range("D44").Value = "'8.097"
$ws.Range("E44").Value = "  -3.65%  "

# Row 45
$ws.Range("D45").Value = "'1.012"
$ws.Range("E45").Value = "  +0.98%  "

# Row 46
$ws.Range("D46").Value = "'0.4588"
$ws.Range("E46").Value = "  -4.50%  "

# Row 47
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Value = "'102.55"
$ws.Range("E47").Value = "  -0.35%  "

# Row 48
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'9.897"
$ws.Range("E48").Value = "  -3.70%  "

# Row 49
$ws.Range("D49").Value = "'1.583"
$ws.Range("E49").Value = "  -4.12%  "

# Row 50
$ws.Range("D50").Value = "'0.06006"
$ws.Range("E50").Value = "  -3.40%  "

# Row 51
$ws.Range("D51").Value = "'63.78"
$ws.Range("E51").Value = "  -2.32%  "
